$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29 (DRAIAM070) keeps its TCID but gets the updated IPAIAM customer
# care script content (OPQA-5154 gains a second Jira id, and the
# description picks up a second "Call us" verification clause).
$ws.Range("A29").Value = 'DRAIAM070'
$ws.Range("B29").Value = 'OPQA-5154||OPQA-5230'
$ws.Range("C29").Value = 'Verify that a "Call us" section is present in customer care contact page with customer care contact details (region, phone numbers, hours of operation, language supported||Ensure that the page has "Support Request" and "Call us" sections.'

# New row 29 content wraps onto two lines, so enable wrap + grow the row.
$ws.Range("C29").WrapText = $true
$ws.Rows.Item(29).RowHeight = 30

# Row 30 (DRAIAM071) is untouched content-wise; values re-asserted here for
# robustness in case the row ever drifts.
$ws.Range("A30").Value = 'DRAIAM071'
$ws.Range("B30").Value = 'OPQA-5168 || OPQA-5227'
$ws.Range("C30").Value = "Verify that the web form provided to user should be application specific with following required fields`n1.Name 2.Organization 3.Contact details (email, telephone) 4.Issue Category 5.Country`n6.Description of issue ( a free form text box where a user can describe why they are contacting support) ||`nVerify that the web form provided to user will be application specific with following required fields (Name, Organization, email, telephone, Issue Category, Country, Description of issue)"

# Move the active selection to B29, matching the updated cursor position.
$ws.Range("B29").Select()
